$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix trailing wording in the warranty-period answer (D4): remove the
# stray space before "합니다." at the very end of the text.
$ws.Range("D4").Value = "고객이 차량의 보증 기간에 대해 묻는 경우 아래와 같이 안내합니다.`n<br>`n<br>1. 차량의 보증 기간은 5년 또는 10만km 중 선도래 우선 기준 (XC60 중국 생산분은 7년 또는 14만km 중 선도래 기준)`n<br>2. 잔여 보증 기간에 대해 묻는 경우 서비스센터로 문의가 필요함을 안내합니다."

# Give the NUGU Auto / TMAP Auto / FLO rows their own Category instead of
# the generic "All", and fix the "Flo" capitalization to "FLO".
$ws.Range("E11").Value = "NUGU Auto"
$ws.Range("E12").Value = "TMAP Auto"
$ws.Range("E14").Value = "TMAP Auto"
$ws.Range("A15").Value = "FLO앱이 동작하지 않는 경우"
$ws.Range("E15").Value = "FLO"

# Update the active selection/scroll position to match the author's saved
# view (scrolled to top, D5 selected).
$ws.Range("D5").Select()
